# Handled the customer not eligible for the agreement.
# Update the CO_Number column (AB) on the "Input" sheet:
#  - Replace the existing CO numbers for rows 2-6 with new values (entered as text)
#  - Clear the CO numbers for rows 7-9 (customer not eligible for the agreement)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Input")

$ws.Range("AB2").Value = "3013697318"
$ws.Range("AB3").Value = "3013697320"
$ws.Range("AB4").Value = "3013697321"
$ws.Range("AB5").Value = "3013697322"
$ws.Range("AB6").Value = "3013697323"

$ws.Range("AB7").ClearContents()
$ws.Range("AB8").ClearContents()
$ws.Range("AB9").ClearContents()

# Reflect the user's resulting selection/scroll position on the sheet.
$ws.Activate()
$ws.Range("AB2:AB10").Select()
$excel.ActiveWindow.ScrollColumn = 17
